$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.04464623277901936
$ws.Range("D2").Value = 0.05120756899277268
$ws.Range("E2").Value = 0.1231914243850269
$ws.Range("F2").Value = 3.261719853951675
$ws.Range("G2").Value = 0.00256223626207244
$ws.Range("I2").Value = 1.956199376791929
$ws.Range("J2").Value = 0.22557843212347
$ws.Range("K2").Value = 2.756743983492072

$ws.Range("C3").Value = 0.04488352438674603
$ws.Range("D3").Value = 0.04998413789638079
$ws.Range("E3").Value = 0.1206574810010821
$ws.Range("F3").Value = 3.239797117299972
$ws.Range("G3").Value = 0.002568249987252136
$ws.Range("I3").Value = 1.940452123511534
$ws.Range("J3").Value = 0.2211319605938726
$ws.Range("K3").Value = 2.585126526994713

$ws.Range("C4").Value = 0.04503551992229582
$ws.Range("D4").Value = 0.04924358304448617
$ws.Range("E4").Value = 0.1191657467624339
$ws.Range("F4").Value = 3.228481138575887
$ws.Range("G4").Value = 0.002572132344621392
$ws.Range("I4").Value = 1.932140555824589
$ws.Range("J4").Value = 0.2185449907730614
$ws.Range("K4").Value = 2.481170118028786

$ws.Range("C5").Value = 0.04509905269896564
$ws.Range("D5").Value = 0.04894453693743372
$ws.Range("E5").Value = 0.118573948526393
$ws.Range("F5").Value = 3.22440582127534
$ws.Range("G5").Value = 0.002573762365045627
$ws.Range("I5").Value = 1.929092505581821
$ws.Range("J5").Value = 0.2175265770613493
$ws.Range("K5").Value = 2.439161996027906

$ws.Range("C6").Value = 0.04510969879759585
$ws.Range("D6").Value = 0.04889504763417563
$ws.Range("E6").Value = 0.1184766522604725
$ws.Range("F6").Value = 3.223761398291089
$ws.Range("G6").Value = 0.002574035928300818
$ws.Range("I6").Value = 1.928606785020591
$ws.Range("J6").Value = 0.2173596263405173
$ws.Range("K6").Value = 2.432207977492112

$ws.Range("C7").Value = 0.04503637028257579
$ws.Range("D7").Value = 0.04923953884367904
$ws.Range("E7").Value = 0.1191577004280902
$ws.Range("F7").Value = 3.228424011308661
$ws.Range("G7").Value = 0.002572154133496776
$ws.Range("I7").Value = 1.932098079286021
$ws.Range("J7").Value = 0.2185311114169082
$ws.Range("K7").Value = 2.480602146437718

$ws.Range("C8").Value = 0.0447267496773236
$ws.Range("D8").Value = 0.05078355820105429
$ws.Range("E8").Value = 0.1223043951222245
$ws.Range("F8").Value = 3.253713786064495
$ws.Range("G8").Value = 0.00256427048473413
$ws.Range("I8").Value = 1.95048660152726
$ws.Range("J8").Value = 0.2240154442100817
$ws.Range("K8").Value = 2.697275136517362

$ws.Range("C9").Value = 0.04416912428142794
$ws.Range("D9").Value = 0.05389331329784142
$ws.Range("E9").Value = 0.1289855866347303
$ws.Range("F9").Value = 3.320478019348684
$ws.Range("G9").Value = 0.002550309503097509
$ws.Range("I9").Value = 1.997425481733799
$ws.Range("J9").Value = 0.2359165469682409
$ws.Range("K9").Value = 3.133511844213615

$ws.Range("C10").Value = 0.04378906433940699
$ws.Range("D10").Value = 0.05622503614063135
$ws.Range("E10").Value = 0.1342087894404642
$ws.Range("F10").Value = 3.38022126288655
$ws.Range("G10").Value = 0.002540954872361504
$ws.Range("I10").Value = 2.038701681819248
$ws.Range("J10").Value = 0.2453746955705896
$ws.Range("K10").Value = 3.461103402948083

$ws.Range("C11").Value = 0.04362248331504759
$ws.Range("D11").Value = 0.05729541376341984
$ws.Range("E11").Value = 0.1366540540424523
$ws.Range("F11").Value = 3.409770929657384
$ws.Range("G11").Value = 0.002536892773279521
$ws.Range("I11").Value = 2.05898851279207
$ws.Range("J11").Value = 0.2498360785733666
$ws.Range("K11").Value = 3.611714384661127

$ws.Range("C12").Value = 0.04356030243461362
$ws.Range("D12").Value = 0.05770207547995909
$ws.Range("E12").Value = 0.1375900172951745
$ws.Range("F12").Value = 3.421305622291186
$ws.Range("G12").Value = 0.002535382183783449
$ws.Range("I12").Value = 2.066890502679328
$ws.Range("J12").Value = 0.2515485929320818
$ws.Range("K12").Value = 3.668977975758992

$ws.Range("C13").Value = 0.0435736543018308
$ws.Range("D13").Value = 0.05761443521131326
$ws.Range("E13").Value = 0.1373879956345476
$ws.Range("F13").Value = 3.418806021060533
$ws.Range("G13").Value = 0.002535706289750229
$ws.Range("I13").Value = 2.065178849114147
$ws.Range("J13").Value = 0.2511787416134723
$ws.Range("K13").Value = 3.656634950791783

$ws.Range("C14").Value = 0.04361734966741437
$ws.Range("D14").Value = 0.05732884367432689
$ws.Range("E14").Value = 0.1367308557389322
$ws.Range("F14").Value = 3.41071296225266
$ws.Range("G14").Value = 0.002536767943232301
$ws.Range("I14").Value = 2.059634193653906
$ws.Range("J14").Value = 0.2499765040632127
$ws.Range("K14").Value = 3.616420862902544

$ws.Range("C15").Value = 0.04364423129211659
$ws.Range("D15").Value = 0.05715408253598753
$ws.Range("E15").Value = 0.1363296414984276
$ws.Range("F15").Value = 3.405800755345098
$ws.Range("G15").Value = 0.002537421831336105
$ws.Range("I15").Value = 2.056266636909982
$ws.Range("J15").Value = 0.2492431126611194
$ws.Range("K15").Value = 3.59181865386887

$ws.Range("C16").Value = 0.04380007711394818
$ws.Range("D16").Value = 0.05615527403794118
$ws.Range("E16").Value = 0.1340503818392946
$ws.Range("F16").Value = 3.378338177155513
$ws.Range("G16").Value = 0.002541224216611404
$ws.Range("I16").Value = 2.037406503067871
$ws.Range("J16").Value = 0.2450863506265648
$ws.Range("K16").Value = 3.451292754399105

$ws.Range("C17").Value = 0.04389729408970311
$ws.Range("D17").Value = 0.0555449719207104
$ws.Range("E17").Value = 0.1326698845233523
$ws.Range("F17").Value = 3.362101005161321
$ws.Range("G17").Value = 0.002543606264587935
$ws.Range("I17").Value = 2.026225072845961
$ws.Range("J17").Value = 0.2425771644266206
$ws.Range("K17").Value = 3.365492858111566

$ws.Range("C18").Value = 0.04395380509738622
$ws.Range("D18").Value = 0.05519485533012158
$ws.Range("E18").Value = 0.131882369999694
$ws.Range("F18").Value = 3.352984865295497
$ws.Range("G18").Value = 0.002544994566427318
$ws.Range("I18").Value = 2.019935785142025
$ws.Range("J18").Value = 0.2411488810941051
$ws.Range("K18").Value = 3.316292356030317

$ws.Range("C19").Value = 0.04397304107687461
$ws.Range("D19").Value = 0.05507647044287722
$ws.Range("E19").Value = 0.1316168481097293
$ws.Range("F19").Value = 3.349936492476274
$ws.Range("G19").Value = 0.002545467754604974
$ws.Range("I19").Value = 2.017830644185949
$ws.Range("J19").Value = 0.240667846420564
$ws.Range("K19").Value = 3.299659495657806

$ws.Range("C20").Value = 0.04388688371820137
$ws.Range("D20").Value = 0.05560984554293924
$ws.Range("E20").Value = 0.1328161666481904
$ws.Range("F20").Value = 3.363806367711533
$ws.Range("G20").Value = 0.002543350807957722
$ws.Range("I20").Value = 2.027400642116902
$ws.Range("J20").Value = 0.2428427241752331
$ws.Range("K20").Value = 3.374610933471047

$ws.Range("C21").Value = 0.04360449091956298
$ws.Range("D21").Value = 0.05741269310206576
$ws.Range("E21").Value = 0.1369236020477871
$ws.Range("F21").Value = 3.413080699794278
$ws.Range("G21").Value = 0.002536455361214695
$ws.Range("I21").Value = 2.061256806517861
$ws.Range("J21").Value = 0.250329001851938
$ws.Range("K21").Value = 3.628226442489051

$ws.Range("C22").Value = 0.04342517198166007
$ws.Range("D22").Value = 0.05859870221127039
$ws.Range("E22").Value = 0.1396663266343836
$ws.Range("F22").Value = 3.447296298174763
$ws.Range("G22").Value = 0.002532109812672962
$ws.Range("I22").Value = 2.084666227551281
$ws.Range("J22").Value = 0.2553563946484587
$ws.Range("K22").Value = 3.795323337707373

$ws.Range("C23").Value = 0.04352040074423513
$ws.Range("D23").Value = 0.05796501668135079
$ws.Range("E23").Value = 0.1381971355226028
$ws.Range("F23").Value = 3.428849439236529
$ws.Range("G23").Value = 0.002534414434744022
$ws.Range("I23").Value = 2.072053947299409
$ws.Range("J23").Value = 0.2526607747046654
$ws.Range("K23").Value = 3.706016820106129

$ws.Range("C24").Value = 0.04389158831618012
$ws.Range("D24").Value = 0.05558051382469387
$ws.Range("E24").Value = 0.1327500133413366
$ws.Range("F24").Value = 3.363034691990407
$ws.Range("G24").Value = 0.002543466241253882
$ws.Range("I24").Value = 2.026868734292407
$ws.Range("J24").Value = 0.2427226201481574
$ws.Range("K24").Value = 3.370488256179954

$ws.Range("C25").Value = 0.04431473388715723
$ws.Range("D25").Value = 0.05304358222103644
$ws.Range("E25").Value = 0.1271231589882049
$ws.Range("F25").Value = 3.300555270540031
$ws.Range("G25").Value = 0.002553927020508039
$ws.Range("I25").Value = 1.983546348395492
$ws.Range("J25").Value = 0.2325726660076839
$ws.Range("K25").Value = 3.014266470208554
